$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.8417533333333332
$ws.Range("H2").Value = 2.52526
$ws.Range("I2").Value = 0.01079423211523897
$ws.Range("J2").Value = 0.01079423211523897
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.02564166666666666
$ws.Range("N2").Value = 0.07692499999999999
$ws.Range("O2").Value = 0.0006780701807970013
$ws.Range("P2").Value = 0.0006780701807970013
$ws.Range("Q2").Value = 0.02158395838888888
$ws.Range("R2").Value = 0.1942556255
$ws.Range("S2").Value = 0.000007319246921944885
$ws.Range("T2").Value = 0.000007319246921944883

$ws.Range("G3").Value = 0.8417533333333332
$ws.Range("H3").Value = 2.52526
$ws.Range("I3").Value = 0.01079423211523897
$ws.Range("J3").Value = 0.01079423211523897
$ws.Range("M3").Value = 0.01112833333333333
$ws.Range("O3").Value = 0.0002942784918545062
$ws.Range("P3").Value = 0.0002942784918545062
$ws.Range("Q3").Value = 0.009367311677777776
$ws.Range("R3").Value = 0.08430580509999999
$ws.Range("S3").Value = 0.0000031765103476
$ws.Range("T3").Value = 0.000003176510347599999

$ws.Range("G4").Value = 0.8417533333333332
$ws.Range("H4").Value = 2.52526
$ws.Range("I4").Value = 0.01079423211523897
$ws.Range("J4").Value = 0.01079423211523897
$ws.Range("O4").Value = 0.0009221225577320236
$ws.Range("P4").Value = 0.0009221225577320235
$ws.Range("Q4").Value = 0.02935249990222222
$ws.Range("R4").Value = 0.26417249912
$ws.Range("S4").Value = 0.000009953604926857308
$ws.Range("T4").Value = 0.000009953604926857305

$ws.Range("G5").Value = 0.8417533333333332
$ws.Range("H5").Value = 2.52526
$ws.Range("I5").Value = 0.01079423211523897
$ws.Range("J5").Value = 0.01079423211523897
$ws.Range("M5").Value = 37.74401233333333
$ws.Range("N5").Value = 113.232037
$ws.Range("O5").Value = 0.9981055287696164
$ws.Range("P5").Value = 0.9981055287696164
$ws.Range("Q5").Value = 31.77114819495777
$ws.Range("R5").Value = 285.94033375462
$ws.Range("S5").Value = 0.01077378275304257
$ws.Range("T5").Value = 0.01077378275304256

$ws.Range("G6").Value = 69.05064766666666
$ws.Range("I6").Value = 0.885471656726338
$ws.Range("J6").Value = 0.8854716567263378
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.02564166666666666
$ws.Range("N6").Value = 0.07692499999999999
$ws.Range("O6").Value = 0.0006780701807970013
$ws.Range("P6").Value = 0.0006780701807970013
$ws.Range("Q6").Value = 1.770573690586111
$ws.Range("R6").Value = 15.935163215275
$ws.Range("S6").Value = 0.0006004119263670483
$ws.Range("T6").Value = 0.0006004119263670481

$ws.Range("G7").Value = 69.05064766666666
$ws.Range("I7").Value = 0.885471656726338
$ws.Range("J7").Value = 0.8854716567263378
$ws.Range("M7").Value = 0.01112833333333333
$ws.Range("O7").Value = 0.0002942784918545062
$ws.Range("P7").Value = 0.0002942784918545062
$ws.Range("Q7").Value = 0.7684186241172221
$ws.Range("R7").Value = 6.915767617054999
$ws.Range("S7").Value = 0.0002605752637213378
$ws.Range("T7").Value = 0.0002605752637213377

$ws.Range("G8").Value = 69.05064766666666
$ws.Range("I8").Value = 0.885471656726338
$ws.Range("J8").Value = 0.8854716567263378
$ws.Range("O8").Value = 0.0009221225577320236
$ws.Range("P8").Value = 0.0009221225577320235
$ws.Range("S8").Value = 0.0008165133888997032
$ws.Range("T8").Value = 0.0008165133888997029

$ws.Range("G9").Value = 69.05064766666666
$ws.Range("I9").Value = 0.885471656726338
$ws.Range("J9").Value = 0.8854716567263378
$ws.Range("M9").Value = 37.74401233333333
$ws.Range("N9").Value = 113.232037
$ws.Range("O9").Value = 0.9981055287696164
$ws.Range("P9").Value = 0.9981055287696164
$ws.Range("Q9").Value = 2606.248497155321
$ws.Range("R9").Value = 23456.23647439789
$ws.Range("S9").Value = 0.8837941561473499
$ws.Range("T9").Value = 0.8837941561473497

$ws.Range("G10").Value = 7.697976666666666
$ws.Range("H10").Value = 23.09393
$ws.Range("I10").Value = 0.09871507918910555
$ws.Range("J10").Value = 0.09871507918910553
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.02564166666666666
$ws.Range("N10").Value = 0.07692499999999999
$ws.Range("O10").Value = 0.0006780701807970013
$ws.Range("P10").Value = 0.0006780701807970013
$ws.Range("Q10").Value = 0.1973889516944444
$ws.Range("R10").Value = 1.77650056525
$ws.Range("S10").Value = 0.0000669357515931471
$ws.Range("T10").Value = 0.00006693575159314709

$ws.Range("G11").Value = 7.697976666666666
$ws.Range("H11").Value = 23.09393
$ws.Range("I11").Value = 0.09871507918910555
$ws.Range("J11").Value = 0.09871507918910553
$ws.Range("M11").Value = 0.01112833333333333
$ws.Range("O11").Value = 0.0002942784918545062
$ws.Range("P11").Value = 0.0002942784918545062
$ws.Range("Q11").Value = 0.08566565033888887
$ws.Range("R11").Value = 0.7709908530499999
$ws.Range("S11").Value = 0.00002904972462706813
$ws.Range("T11").Value = 0.00002904972462706813

$ws.Range("G12").Value = 7.697976666666666
$ws.Range("H12").Value = 23.09393
$ws.Range("I12").Value = 0.09871507918910555
$ws.Range("J12").Value = 0.09871507918910553
$ws.Range("O12").Value = 0.0009221225577320236
$ws.Range("P12").Value = 0.0009221225577320235
$ws.Range("Q12").Value = 0.2684335783511111
$ws.Range("R12").Value = 2.41590220516
$ws.Range("S12").Value = 0.00009102740130857727
$ws.Range("T12").Value = 0.00009102740130857724

$ws.Range("G13").Value = 7.697976666666666
$ws.Range("H13").Value = 23.09393
$ws.Range("I13").Value = 0.09871507918910555
$ws.Range("J13").Value = 0.09871507918910553
$ws.Range("M13").Value = 37.74401233333333
$ws.Range("N13").Value = 113.232037
$ws.Range("O13").Value = 0.9981055287696164
$ws.Range("P13").Value = 0.9981055287696164
$ws.Range("Q13").Value = 290.5525262483789
$ws.Range("R13").Value = 2614.97273623541
$ws.Range("S13").Value = 0.09852806631157675
$ws.Range("T13").Value = 0.09852806631157673

$ws.Range("G14").Value = 0.391393
$ws.Range("H14").Value = 1.174179
$ws.Range("I14").Value = 0.005019031969317685
$ws.Range("J14").Value = 0.005019031969317684
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.02564166666666666
$ws.Range("N14").Value = 0.07692499999999999
$ws.Range("O14").Value = 0.0006780701807970013
$ws.Range("P14").Value = 0.0006780701807970013
$ws.Range("Q14").Value = 0.01003596884166667
$ws.Range("R14").Value = 0.090323719575
$ws.Range("S14").Value = 0.000003403255914861173
$ws.Range("T14").Value = 0.000003403255914861172

$ws.Range("G15").Value = 0.391393
$ws.Range("H15").Value = 1.174179
$ws.Range("I15").Value = 0.005019031969317685
$ws.Range("J15").Value = 0.005019031969317684
$ws.Range("M15").Value = 0.01112833333333333
$ws.Range("O15").Value = 0.0002942784918545062
$ws.Range("P15").Value = 0.0002942784918545062
$ws.Range("Q15").Value = 0.004355551768333334
$ws.Range("R15").Value = 0.039199965915
$ws.Range("S15").Value = 0.000001476993158500361
$ws.Range("T15").Value = 0.000001476993158500361

$ws.Range("G16").Value = 0.391393
$ws.Range("H16").Value = 1.174179
$ws.Range("I16").Value = 0.005019031969317685
$ws.Range("J16").Value = 0.005019031969317684
$ws.Range("O16").Value = 0.0009221225577320236
$ws.Range("P16").Value = 0.0009221225577320235
$ws.Range("Q16").Value = 0.01364813483866667
$ws.Range("R16").Value = 0.122833213548
$ws.Range("S16").Value = 0.00000462816259688602
$ws.Range("T16").Value = 0.000004628162596886018

$ws.Range("G17").Value = 0.391393
$ws.Range("H17").Value = 1.174179
$ws.Range("I17").Value = 0.005019031969317685
$ws.Range("J17").Value = 0.005019031969317684
$ws.Range("M17").Value = 37.74401233333333
$ws.Range("N17").Value = 113.232037
$ws.Range("O17").Value = 0.9981055287696164
$ws.Range("P17").Value = 0.9981055287696164
$ws.Range("Q17").Value = 14.77274221918033
$ws.Range("R17").Value = 132.954679972623
$ws.Range("S17").Value = 0.005009523557647438
$ws.Range("T17").Value = 0.005009523557647437
